# Reproduce the manual edits the author made in Excel before their DB-import
# testing session: they widened column A (PLACA) slightly and left the
# selection on M16 instead of A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to a custom width (author resized it, landing close to the
# same width already used by the stray column H from an earlier edit).
$ws.Columns.Item(1).ColumnWidth = 9.5

# Leave the active selection on M16, matching the saved cursor position.
$ws.Range("M16").Select()
